# Apply the "Tracks rem identifiers column" edit:
#  - Rename header "Identifiers" (O1) -> "Internal  House  Name"
#  - Add a new trailing column S with header "Tag /Band"
#  - Give column O (the renamed column) a wider custom width, and bump
#    the sheet's default column width slightly
#  - Move the active selection / view over towards the new column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text changes -----------------------------------------------
$ws.Range("O1").Value = "Internal  House  Name"
$ws.Range("S1").Value = "Tag /Band"

# --- Column widths -------------------------------------------------------
# Column O (the renamed column) gets a wider, explicit width.
$ws.Columns.Item(15).ColumnWidth = 14.6

# Bump the sheet-wide standard/default column width slightly.
$ws.StandardWidth = 10.74

# --- View / selection ------------------------------------------------------
# Scroll the view toward the new columns and move the active selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 11
$ws.Range("O1").Select()
